$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 508; this shifts existing rows 508..557 down to 509..558
$ws.Rows(508).Insert()

# Populate the newly inserted row 508 with the new data record
$ws.Cells.Item(508, 1).Value = 11
$ws.Cells.Item(508, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(508, 3).Value = "Bíobío"
$ws.Cells.Item(508, 4).Value = 45212
$ws.Cells.Item(508, 5).Value = 8
$ws.Cells.Item(508, 6).Value = 100112017
$ws.Cells.Item(508, 7).Value = "Apio"
$ws.Cells.Item(508, 8).Value = "Americana (o)"
$ws.Cells.Item(508, 9).Value = "Primera"
$ws.Cells.Item(508, 10).Value = 250
$ws.Cells.Item(508, 11).Value = 8000
$ws.Cells.Item(508, 12).Value = 9000
$ws.Cells.Item(508, 13).Value = 8400
$ws.Cells.Item(508, 14).Value = "$/docena de matas"
$ws.Cells.Item(508, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(508, 16).Value = 1400
$ws.Cells.Item(508, 17).Value = 6
$ws.Cells.Item(508, 18).Value = "Hortaliza"
